# "Generate Report for Handoff"
#
# The e315489c-191c-4f5d-87f0-877eecaa0c2d.md row (row 6 on every sheet) has
# just been handed off for localization, so its status flips from
# "In Translation" to "Ready for handoff" and the various handoff timestamps
# advance to reflect the new handoff run.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E6").Value = "Ready for handoff"          # zh-cn status
$ws.Range("F6").Value = "Ready for handoff"          # de-de status
$ws.Range("G6").Value = "2016-08-31 07:24:40"        # Latest HO Xliff Generate Date

# ---- zh-cn detail sheet ------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C6").Value = "Ready for handoff"          # Status
$ws.Range("H6").Value = "2016-08-31 07:24:29"        # Latest Handoff Datetime

# ---- de-de detail sheet ------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C6").Value = "Ready for handoff"          # Status
$ws.Range("H6").Value = "2016-08-31 07:24:40"        # Latest Handoff Datetime
